$wb = $excel.ActiveWorkbook

# Mapping of sheet number -> list of (row, expectedOldValue, newValue) for column F ("想去人数")
$changes = @{
    "展览" = @(@{Row=6; Old=1601; New=1603}, @{Row=8; Old=40; New=41}, @{Row=9; Old=732; New=735}, @{Row=10; Old=2685; New=2690}, @{Row=11; Old=2685; New=2690}, @{Row=12; Old=17; New=18}, @{Row=13; Old=1777; New=1787}, @{Row=15; Old=277; New=284}, @{Row=16; Old=698; New=699}, @{Row=17; Old=5097; New=5126}, @{Row=19; Old=74; New=77}, @{Row=20; Old=689; New=690}, @{Row=22; Old=865; New=870}, @{Row=24; Old=71; New=72}, @{Row=25; Old=41; New=43}, @{Row=26; Old=2431; New=2439}, @{Row=28; Old=371; New=372}, @{Row=29; Old=19; New=20}, @{Row=31; Old=483; New=486}, @{Row=32; Old=1304; New=1307}, @{Row=33; Old=804; New=805}, @{Row=34; Old=7; New=8}, @{Row=35; Old=67; New=68}, @{Row=36; Old=23; New=25}, @{Row=38; Old=1453; New=1458}, @{Row=39; Old=19; New=20}, @{Row=40; Old=1402; New=1409})
    "演出" = @(@{Row=9; Old=114; New=115}, @{Row=11; Old=139; New=148}, @{Row=13; Old=78; New=79}, @{Row=17; Old=333; New=334}, @{Row=18; Old=260; New=261}, @{Row=19; Old=517; New=518})
    "本地生活" = @(@{Row=3; Old=886; New=892}, @{Row=4; Old=250; New=251}, @{Row=6; Old=39; New=40}, @{Row=7; Old=62; New=63})
    "全部类型" = @(@{Row=4; Old=886; New=892}, @{Row=5; Old=250; New=251}, @{Row=8; Old=39; New=40}, @{Row=9; Old=62; New=63}, @{Row=13; Old=1601; New=1603}, @{Row=16; Old=40; New=41}, @{Row=17; Old=2685; New=2690}, @{Row=19; Old=17; New=18}, @{Row=20; Old=1777; New=1787}, @{Row=21; Old=139; New=148}, @{Row=23; Old=277; New=284}, @{Row=24; Old=698; New=699}, @{Row=25; Old=5097; New=5126}, @{Row=27; Old=74; New=77}, @{Row=28; Old=689; New=690}, @{Row=30; Old=865; New=870}, @{Row=32; Old=71; New=72}, @{Row=34; Old=41; New=43}, @{Row=35; Old=2431; New=2439}, @{Row=37; Old=371; New=372}, @{Row=39; Old=483; New=487}, @{Row=40; Old=1304; New=1307}, @{Row=42; Old=260; New=261}, @{Row=43; Old=517; New=518}, @{Row=44; Old=804; New=805}, @{Row=45; Old=7; New=8}, @{Row=46; Old=67; New=68}, @{Row=47; Old=23; New=25}, @{Row=49; Old=1402; New=1409})
}

$totalApplied = 0
$mismatchCount = 0

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($chg in $changes[$sheetName]) {
        $cell = $ws.Cells.Item($chg.Row, 6)
        $current = $cell.Value2
        if ($current -ne $chg.Old) {
            $mismatchCount++
            Write-Host "WARNING: $sheetName row $($chg.Row) expected $($chg.Old) but found $current"
        }
        $cell.Value = $chg.New
        $totalApplied++
    }
}

Write-Host "Done applying $totalApplied changes ($mismatchCount mismatches)."